$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s)*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Docente(s) Responsável(eis)' paragraph"
}

# Insert a new empty paragraph right after it.
$target.Range.InsertParagraphAfter()

# Re-locate that new (currently empty) paragraph: it now sits right after
# the heading paragraph, and right before "Programa resumido".
$newPara = $null
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $newPara = $p
        break
    }
    if ($p.Range.Text -like "Docente(s)*") {
        $found = $true
    }
}

# Fill it in with the "List Bullet" style and the two professors, the first
# one followed by a manual line break, second on its own run - authored as
# raw WordprocessingML so the run/style structure comes out exactly right.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>3380737 - Flávio Teixeira da Silva</w:t><w:br/></w:r><w:r><w:t>8853480 - Tatiane da Franca Silva</w:t></w:r></w:p>'
$newPara.Range.InsertXML($xml)
